# Updated cryptos list with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '30.313.50'
$ws.Range('E2').Value = '  -0.48%  '

$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '1.883.22'
$ws.Range('E3').Value = '  -1.35%  '

$ws.Range('D4').NumberFormat = "@"
$ws.Range('D4').Value = '1.000'
$ws.Range('E4').Value = '  -0.06%  '

$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '238.05'
$ws.Range('E5').Value = '  -0.23%  '

$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '1.001'
$ws.Range('E6').Value = '  +0.04%  '

$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.4693'
$ws.Range('E7').Value = '  -1.59%  '

$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.2841'
$ws.Range('E8').Value = '  +0.57%  '

$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.06597'
$ws.Range('E9').Value = '  -1.31%  '

$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '20.82'
$ws.Range('E10').Value = '  +11.92%  '

$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.07788'
$ws.Range('E11').Value = '  +1.39%  '

$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '98.02'
$ws.Range('E12').Value = '  -2.33%  '

$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '1.886.98'
$ws.Range('E13').Value = '  -1.20%  '

$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '5.101'
$ws.Range('E14').Value = '  -1.65%  '

$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '0.6771'
$ws.Range('E15').Value = '  +1.68%  '

$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '283.80'
$ws.Range('E16').Value = '  +11.44%  '

$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '30.312.90'
$ws.Range('E17').Value = '  -0.58%  '

$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '0.9991'
$ws.Range('E18').Value = '  -0.15%  '

$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '12.65'
$ws.Range('E19').Value = '  +0.32%  '

$ws.Range('B20').Value = 'Uniswap'
$ws.Range('C20').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '5.407'
$ws.Range('E20').Value = '  +0.84%  '

$ws.Range('B21').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C21').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '2.123.02'
$ws.Range('E21').Value = '  -1.56%  '

$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '0.000007288'
$ws.Range('E22').Value = '  -2.15%  '

$ws.Range('E23').Value = '  +0.00%  '

$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '6.183'
$ws.Range('E24').Value = '  -1.43%  '

$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '9.421'
$ws.Range('E25').Value = '  +1.25%  '

$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '168.18'
$ws.Range('E26').Value = '  +0.37%  '

$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '19.27'
$ws.Range('E27').Value = '  +1.23%  '

$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '1.997'
$ws.Range('E28').Value = '  -2.58%  '

$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '1.378'
$ws.Range('E29').Value = '  -0.15%  '

$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '0.09734'
$ws.Range('E30').Value = '  -2.46%  '

$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '4.410'
$ws.Range('E31').Value = '  -7.68%  '

$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '1.488'
$ws.Range('E32').Value = '  -1.23%  '

$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '4.138'
$ws.Range('E33').Value = '  -2.81%  '

$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '0.04681'
$ws.Range('E34').Value = '  -0.24%  '

$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '0.7086'
$ws.Range('E35').Value = '  -1.93%  '

$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '1.100'
$ws.Range('E36').Value = '  -0.08%  '

$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '2.714'
$ws.Range('E37').Value = '  +0.43%  '

$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '0.01877'
$ws.Range('E38').Value = '  -1.58%  '

$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '6.710'
$ws.Range('E39').Value = '  +7.27%  '

$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '2.527'
$ws.Range('E40').Value = '  -3.09%  '

$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '72.29'
$ws.Range('E41').Value = '  -3.26%  '

$ws.Range('B42').Value = 'RenderToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '1.975'
$ws.Range('E42').Value = '  +1.17%  '

$ws.Range('B43').Value = 'TrustWalletToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '0.8656'
$ws.Range('E43').Value = '  +0.83%  '

$ws.Range('E44').Value = '  +0.10%  '

$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '103.59'
$ws.Range('E45').Value = '  -1.75%  '

$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '0.4196'
$ws.Range('E46').Value = '  -0.65%  '

$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '992.25'
$ws.Range('E47').Value = '  +8.71%  '

$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '7.295'
$ws.Range('E48').Value = '  -0.50%  '

$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '9.244'
$ws.Range('E49').Value = '  +6.01%  '

$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '34.00'
$ws.Range('E50').Value = '  -1.84%  '

$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '0.1151'
$ws.Range('E51').Value = '  -3.74%  '
